$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("13/04/2022", "16:39", "17:32", 14, 10, 4, 1, 7, "EURUSD", 15.17, "LOSS"),
    @("14/04/2022", "10:59", "12:08", 5, 2, 3, 2, 2, "EURUSD", 0.9199999999999999, "LOSS"),
    @("14/04/2022", "12:08", "12:11", 3, 0, 3, 3, 0, "EURUSD", -14, "LOSS"),
    @("14/04/2022", "12:11", "12:18", 7, 3, 4, 2, 2, "EURUSD", -1.08, "LOSS")
)

$startRow = 8
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}
